$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = 44964
$ws.Cells.Item(2, 10).Value = 300
$ws.Cells.Item(2, 11).Value = 20000
$ws.Cells.Item(2, 12).Value = 21000
$ws.Cells.Item(2, 13).Value = 20500
$ws.Cells.Item(2, 16).Value = 1139
$ws.Cells.Item(3, 4).Value = 45229
$ws.Cells.Item(3, 10).Value = 460
$ws.Cells.Item(3, 11).Value = 16000
$ws.Cells.Item(3, 12).Value = 17000
$ws.Cells.Item(3, 13).Value = 16500
$ws.Cells.Item(3, 16).Value = 917
$ws.Cells.Item(4, 4).Value = 44984
$ws.Cells.Item(4, 10).Value = 200
$ws.Cells.Item(4, 11).Value = 17000
$ws.Cells.Item(4, 12).Value = 18000
$ws.Cells.Item(4, 13).Value = 17500
$ws.Cells.Item(4, 16).Value = 972
$ws.Cells.Item(5, 4).Value = 44957
$ws.Cells.Item(5, 10).Value = 400
$ws.Cells.Item(5, 11).Value = 21000
$ws.Cells.Item(5, 12).Value = 22000
$ws.Cells.Item(5, 13).Value = 21500
$ws.Cells.Item(5, 16).Value = 1194
$ws.Cells.Item(6, 4).Value = 45177
$ws.Cells.Item(6, 10).Value = 540
$ws.Cells.Item(6, 11).Value = 16000
$ws.Cells.Item(6, 12).Value = 17000
$ws.Cells.Item(6, 13).Value = 16500
$ws.Cells.Item(6, 16).Value = 917
$ws.Cells.Item(7, 4).Value = 45117
$ws.Cells.Item(7, 10).Value = 300
$ws.Cells.Item(7, 11).Value = 17000
$ws.Cells.Item(7, 12).Value = 18000
$ws.Cells.Item(7, 13).Value = 17500
$ws.Cells.Item(7, 16).Value = 972
$ws.Cells.Item(8, 4).Value = 45180
$ws.Cells.Item(8, 10).Value = 400
$ws.Cells.Item(8, 11).Value = 16500
$ws.Cells.Item(8, 12).Value = 17000
$ws.Cells.Item(8, 13).Value = 16750
$ws.Cells.Item(8, 16).Value = 931
$ws.Cells.Item(9, 4).Value = 45154
$ws.Cells.Item(9, 10).Value = 500
$ws.Cells.Item(9, 11).Value = 16500
$ws.Cells.Item(9, 12).Value = 17000
$ws.Cells.Item(9, 13).Value = 16750
$ws.Cells.Item(9, 16).Value = 931
$ws.Cells.Item(10, 4).Value = 45222
$ws.Cells.Item(10, 10).Value = 300
$ws.Cells.Item(10, 11).Value = 16000
$ws.Cells.Item(10, 12).Value = 17000
$ws.Cells.Item(10, 13).Value = 16500
$ws.Cells.Item(10, 16).Value = 917
$ws.Cells.Item(12, 4).Value = 44557
$ws.Cells.Item(12, 10).Value = 400
$ws.Cells.Item(12, 11).Value = 13000
$ws.Cells.Item(12, 12).Value = 14000
$ws.Cells.Item(12, 13).Value = 13500
$ws.Cells.Item(12, 16).Value = 750
$ws.Cells.Item(13, 4).Value = 44960
$ws.Cells.Item(13, 10).Value = 400
$ws.Cells.Item(13, 11).Value = 19500
$ws.Cells.Item(13, 12).Value = 20000
$ws.Cells.Item(13, 13).Value = 19750
$ws.Cells.Item(13, 16).Value = 1097
$ws.Cells.Item(14, 4).Value = 45215
$ws.Cells.Item(14, 10).Value = 400
$ws.Cells.Item(14, 11).Value = 16000
$ws.Cells.Item(14, 12).Value = 17000
$ws.Cells.Item(14, 13).Value = 16500
$ws.Cells.Item(14, 16).Value = 917
$ws.Cells.Item(15, 4).Value = 44977
$ws.Cells.Item(15, 10).Value = 400
$ws.Cells.Item(15, 11).Value = 16500
$ws.Cells.Item(15, 12).Value = 17000
$ws.Cells.Item(15, 13).Value = 16750
$ws.Cells.Item(15, 16).Value = 931
$ws.Cells.Item(16, 4).Value = 44998
$ws.Cells.Item(16, 10).Value = 320
$ws.Cells.Item(16, 11).Value = 17000
$ws.Cells.Item(16, 12).Value = 18000
$ws.Cells.Item(16, 13).Value = 17500
$ws.Cells.Item(16, 16).Value = 972
$ws.Cells.Item(17, 4).Value = 45194
$ws.Cells.Item(17, 10).Value = 400
$ws.Cells.Item(17, 11).Value = 16500
$ws.Cells.Item(17, 12).Value = 17000
$ws.Cells.Item(17, 13).Value = 16750
$ws.Cells.Item(17, 16).Value = 931
$ws.Cells.Item(18, 4).Value = 44547
$ws.Cells.Item(18, 10).Value = 200
$ws.Cells.Item(18, 11).Value = 13000
$ws.Cells.Item(18, 12).Value = 14000
$ws.Cells.Item(18, 13).Value = 13500
$ws.Cells.Item(18, 16).Value = 750
$ws.Cells.Item(19, 4).Value = 45159
$ws.Cells.Item(19, 10).Value = 400
$ws.Cells.Item(19, 11).Value = 16000
$ws.Cells.Item(19, 12).Value = 17000
$ws.Cells.Item(19, 13).Value = 16500
$ws.Cells.Item(19, 16).Value = 917
$ws.Cells.Item(20, 4).Value = 45230
$ws.Cells.Item(20, 10).Value = 360
$ws.Cells.Item(20, 11).Value = 16000
$ws.Cells.Item(20, 12).Value = 17000
$ws.Cells.Item(20, 13).Value = 16500
$ws.Cells.Item(20, 16).Value = 917
$ws.Cells.Item(22, 4).Value = 45068
$ws.Cells.Item(22, 10).Value = 400
$ws.Cells.Item(22, 11).Value = 16000
$ws.Cells.Item(22, 12).Value = 17000
$ws.Cells.Item(22, 13).Value = 16500
$ws.Cells.Item(22, 16).Value = 917
$ws.Cells.Item(23, 4).Value = 45166
$ws.Cells.Item(23, 10).Value = 200
$ws.Cells.Item(23, 11).Value = 16000
$ws.Cells.Item(23, 12).Value = 17000
$ws.Cells.Item(23, 13).Value = 16500
$ws.Cells.Item(23, 16).Value = 917
$ws.Cells.Item(24, 4).Value = 45152
$ws.Cells.Item(24, 10).Value = 500
$ws.Cells.Item(24, 11).Value = 16000
$ws.Cells.Item(24, 12).Value = 17000
$ws.Cells.Item(24, 13).Value = 16500
$ws.Cells.Item(24, 16).Value = 917
$ws.Cells.Item(25, 4).Value = 45005
$ws.Cells.Item(25, 10).Value = 200
$ws.Cells.Item(25, 11).Value = 17000
$ws.Cells.Item(25, 12).Value = 18000
$ws.Cells.Item(25, 13).Value = 17500
$ws.Cells.Item(25, 16).Value = 972
